# db_defaultEnergyAssets.xlsx edit
# Commit message: "No more default assets - Config is now the only way to
# add energy assets. Also for households and industry"
#
# Concrete data-level effects (per the canonical OOXML diff):
#  1. The shared string "OTHER_ELECTRICITY_CONSUMPTION" is retired and every
#     cell that used it now uses a new string "ELECTRICITY_CONSUMPTION_PROFILE"
#     instead (sheet "consumptionAssets", cells E2/E4/E5).
#  2. The "storageAssets" sheet gains values in the previously-mostly-empty
#     column P ("vehicle_scaling") for rows 2-16 (0 for regular storage
#     assets, 1 for the two electric-vehicle rows EV/EHGV).
#  3. Minor view/selection bookkeeping: the active/selected sheet moves from
#     "conversionAssets" to "storageAssets", and the remembered selections on
#     a couple of sheets change.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) consumptionAssets: replace OTHER_ELECTRICITY_CONSUMPTION references
# ---------------------------------------------------------------------
$wsConsumption = $wb.Worksheets.Item("consumptionAssets")

$wsConsumption.Range("E2").Value = "ELECTRICITY_CONSUMPTION_PROFILE"
$wsConsumption.Range("E4").Value = "ELECTRICITY_CONSUMPTION_PROFILE"
$wsConsumption.Range("E5").Value = "ELECTRICITY_CONSUMPTION_PROFILE"

$wsConsumption.Range("G7").Select()

# ---------------------------------------------------------------------
# 2) storageAssets: populate the "vehicle_scaling" column (P)
# ---------------------------------------------------------------------
$wsStorage = $wb.Worksheets.Item("storageAssets")

$wsStorage.Range("P2").Value = 0
$wsStorage.Range("P3").Value = 0
$wsStorage.Range("P4").Value = 0
$wsStorage.Range("P5").Value = 0
$wsStorage.Range("P6").Value = 0
$wsStorage.Range("P7").Value = 0
$wsStorage.Range("P8").Value = 0
$wsStorage.Range("P9").Value = 0
$wsStorage.Range("P10").Value = 1
$wsStorage.Range("P11").Value = 0
$wsStorage.Range("P12").Value = 0
$wsStorage.Range("P13").Value = 0
$wsStorage.Range("P14").Value = 0
$wsStorage.Range("P15").Value = 1
$wsStorage.Range("P16").Value = 0
# P17 already has a value (50) and is unchanged.

# ---------------------------------------------------------------------
# 3) View bookkeeping: storageAssets becomes the active/selected sheet
#    (conversionAssets previously had tabSelected="1"); selections move.
# ---------------------------------------------------------------------
$wsConversion = $wb.Worksheets.Item("conversionAssets")
$wsConversion.Range("K15").Select()

$wsStorage.Select()
$wsStorage.Range("P16").Select()
$wsStorage.Application.ActiveWindow.ScrollColumn = 3

$wb.Worksheets.Item("consumptionAssets").Range("E28").Select() | Out-Null
$wsConsumption.Select()
$wsConsumption.Range("E28").Select()

$wsStorage.Select()
